# Financial Performance Data - source data refresh
# Updates the YOY expense/profitability figures in Sheet1 and moves the
# active selection/view to match the author's last on-screen state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: TOTAL_EXPENSES -------------------------------------------------
$ws.Range("E2").Value = 729392
$ws.Range("F2").Value = 779143
$ws.Range("G2").Value = 0.068000000000000005
$ws.Range("H2").Value = 7

# --- Row 5: VISIT_COUNT -----------------------------------------------------
$ws.Range("E5").Value = 198222
$ws.Range("F5").Value = 192148
$ws.Range("G5").Value = -0.031
$ws.Range("H5").Value = -3

# --- Row 12: growth-rate sign fix ------------------------------------------
$ws.Range("G12").Value = 0.97499999999999998
$ws.Range("H12").Value = 98

# --- Row 13 ------------------------------------------------------------
$ws.Range("E13").Value = 566

# --- Rows 21-48: per-period figures ----------------------------------------
$ws.Range("F21").Value = 122566

$ws.Range("F23").Formula = "=F22-F21"

$ws.Range("F24").Value = -18443
$ws.Range("F27").Value = 126280
$ws.Range("F30").Value = -19358
$ws.Range("F33").Value = 122529
$ws.Range("F34").Value = 121642
$ws.Range("F36").Value = -20246
$ws.Range("F39").Value = 129719
$ws.Range("F40").Value = 128779
$ws.Range("F42").Value = -21186
$ws.Range("F48").Value = -22040

# --- View state: scroll + selection -----------------------------------------
$excel.ActiveWindow.ScrollRow = 21
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("E13:F48").Select()
